# Applies the changes described by the diff:
#  1. "Youtube tutorials:" paragraph - drop the spellcheck proofErr wrapper
#     around "Youtube" and merge the two runs into a single run.
#  2. "25/9/2024 ..." paragraph - drop the spellcheck proofErr wrappers
#     around "figma" and "youtube" and merge all runs into a single run.
#  3. "30/9/2024 ..." paragraph - append a new run ", created user profile
#     page" and add a new empty paragraph right after it (before sectPr).

function Set-ParagraphXml($para, $bodyInnerXml) {
    $r = $para.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

function Find-ParagraphContaining($doc, $needle) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

$d = $word.ActiveDocument

# 1) "Youtube tutorials:" - remove proofErr, merge into one run
$pYoutube = Find-ParagraphContaining $d "Youtube tutorials:"
Set-ParagraphXml $pYoutube '<w:p><w:r><w:t>Youtube tutorials:</w:t></w:r></w:p>'

# 2) "25/9/2024 – watched figma tutorials on youtube" - remove proofErr, merge into one run
$pReport1 = Find-ParagraphContaining $d "25/9/2024"
Set-ParagraphXml $pReport1 '<w:p><w:r><w:t>25/9/2024 – watched figma tutorials on youtube</w:t></w:r></w:p>'

# 3) "30/9/2024 – finished article page" - add new run with the extra
#    sentence, then a new empty paragraph right after it.
$pReport4 = Find-ParagraphContaining $d "30/9/2024"
Set-ParagraphXml $pReport4 '<w:p><w:r><w:t>30/9/2024 – finished article page</w:t></w:r><w:r><w:t>, created user profile page</w:t></w:r></w:p><w:p/>'
